$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("Company Name", "Company Number", "Incorporation Date", "Status", "Source", "Date Downloaded", "Time Discovered")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$headerRange = $ws.Range("A1:G1")

$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Borders.Item(7).LineStyle = 1    # xlEdgeLeft, xlContinuous
    $cell.Borders.Item(7).Weight = 2       # xlThin
    $cell.Borders.Item(8).LineStyle = 1    # xlEdgeTop
    $cell.Borders.Item(8).Weight = 2
    $cell.Borders.Item(9).LineStyle = 1    # xlEdgeBottom
    $cell.Borders.Item(9).Weight = 2
    $cell.Borders.Item(10).LineStyle = 1   # xlEdgeRight
    $cell.Borders.Item(10).Weight = 2
}
